$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.699.44"
$ws.Range("E2").Value = "  -6.77%  "

# Row 3
$ws.Range("D3").Value = "2.542.84"
$ws.Range("E3").Value = "  -4.87%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.10"
$ws.Range("E5").Value = "  -4.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.03"
$ws.Range("E6").Value = "  -4.57%  "

# Row 7
$ws.Range("E7").Value = "  -4.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -6.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.16"
$ws.Range("E10").Value = "  -5.69%  "

# Row 11
$ws.Range("E11").Value = "  -5.28%  "

# Row 12
$ws.Range("E12").Value = "  -5.66%  "

# Row 13
$ws.Range("E13").Value = "  +6.66%  "

# Row 14
$ws.Range("D14").Value = "2.935.76"
$ws.Range("E14").Value = "  -5.09%  "

# Row 15
$ws.Range("D15").Value = "2.512.77"
$ws.Range("E15").Value = "  -6.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.875"
$ws.Range("E16").Value = "  -6.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("E17").Value = "  -6.44%  "

# Row 18
$ws.Range("D18").Value = "42.717.43"
$ws.Range("E18").Value = "  -6.76%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0979"
$ws.Range("E19").Value = "  -4.07%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  -4.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.70"
$ws.Range("E22").Value = "  -4.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.04"
$ws.Range("E23").Value = "  -10.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  -5.57%  "

# Row 25
$ws.Range("E25").Value = "  -5.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.96"
$ws.Range("E26").Value = "  -7.69%  "

# Row 27
$ws.Range("E27").Value = "  +0.27%  "

# Row 28
$ws.Range("E28").Value = "  -4.14%  "

# Row 29
$ws.Range("E29").Value = "  -3.61%  "

# Row 30
$ws.Range("E30").Value = "  -2.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  -2.84%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.04"
$ws.Range("E32").Value = "  -2.08%  "

# Row 33
$ws.Range("E33").Value = "  -2.20%  "

# Row 34
$ws.Range("E34").Value = "  -8.96%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.36"
$ws.Range("E35").Value = "  -10.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("E36").Value = "  -6.10%  "

# Row 37
$ws.Range("E37").Value = "  -6.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.17"
$ws.Range("E38").Value = "  +5.50%  "

# Row 39
$ws.Range("E39").Value = "  -4.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.34"
$ws.Range("E40").Value = "  -10.76%  "

# Row 41
$ws.Range("E41").Value = "  -5.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -6.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44
$ws.Range("D44").Value = "2.081.75"
$ws.Range("E44").Value = "  -3.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.68"
$ws.Range("E46").Value = "  +5.10%  "

# Row 47
$ws.Range("E47").Value = "  -3.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.33"
$ws.Range("E48").Value = "  -10.68%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.00"
$ws.Range("E49").Value = "  -6.28%  "

# Row 50
$ws.Range("D50").Value = "2.788.74"
$ws.Range("E50").Value = "  -5.16%  "

# Row 51
$ws.Range("E51").Value = "  -3.49%  "
